$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Farmers Database" - add a new "Black pepper" crop column (I)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New header cell I1, formatted like the other header cells (copy format from H1)
$ws1.Range("I1").Value = "Black pepper"
$ws1.Range("H1").Copy()
$ws1.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing quantities
$ws1.Range("E2").Value = 1      # Arvi - Rice quantity 5 -> 1
$ws1.Range("E4").Value = 1      # Shan - Rice quantity 0 -> 1
$ws1.Range("F4").Value = 400    # Shan - Cashewnuts 500 -> 400

# New "Black pepper" column values for existing rows
$ws1.Range("I2").Value = 0
$ws1.Range("I3").Value = 0
$ws1.Range("I4").Value = 0

# ---------------------------------------------------------------------
# Sheet 2: "Farmers Log" - add a new log entry (row 4)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(4, 1).Value = "2017-10-02 17:32:45"
$ws2.Cells.Item(4, 2).Value = "Shan"
$ws2.Cells.Item(4, 3).Value = 12343
$ws2.Cells.Item(4, 4).Value = "Banana"
$ws2.Cells.Item(4, 5).Value = 5

# ---------------------------------------------------------------------
# Sheet 3: "Companies Log" - add two new log entries (rows 3 and 4)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(3, 1).Value = "2017-10-02 14:29:10"
$ws3.Cells.Item(3, 2).Value = "Parle"
$ws3.Cells.Item(3, 3).Value = "Rice"
$ws3.Cells.Item(3, 4).Value = 100
$ws3.Cells.Item(3, 5).Value = 360000

$ws3.Cells.Item(4, 1).Value = "2017-10-02 17:34:18"
$ws3.Cells.Item(4, 2).Value = "Nestle"
$ws3.Cells.Item(4, 3).Value = "Banana"
$ws3.Cells.Item(4, 4).Value = 8
$ws3.Cells.Item(4, 5).Value = 32800
